$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 2452.25
$ws.Range("I16").Value = 2154.5
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 2154.5
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -1924.5
$ws.Range("N16").Value = -3210
# Row 17
$ws.Range("H17").Value = 2419438.5
$ws.Range("J17").Value = 3002651.5
$ws.Range("L17").Value = 9007954.5
$ws.Range("N17").Value = -9008290.5
# Row 18
$ws.Range("H18").Value = 268.53845
$ws.Range("I18").Value = 268.53845
$ws.Range("K18").Value = 268.53845
$ws.Range("M18").Value = 15.46154999999999
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 33
$ws.Range("H33").Value = 1950119.1
$ws.Range("J33").Value = 4732.5
$ws.Range("L33").Value = 4732.5
$ws.Range("N33").Value = -5190.5
# Row 64
$ws.Range("H64").Value = 5004.625
$ws.Range("I64").Value = 4099.091
$ws.Range("J64").Value = 6996.8
$ws.Range("K64").Value = 4099.091
$ws.Range("L64").Value = 6996.8
$ws.Range("M64").Value = -3851.091
$ws.Range("N64").Value = -7492.8
# Row 67
$ws.Range("H67").Value = 5004.625
$ws.Range("I67").Value = 4099.091
$ws.Range("J67").Value = 6996.8
$ws.Range("K67").Value = 4099.091
$ws.Range("L67").Value = 6996.8
$ws.Range("M67").Value = -3241.091
$ws.Range("N67").Value = -8712.799999999999
# Row 74
$ws.Range("H74").Value = 6561.8184
$ws.Range("I74").Value = 5810.778
$ws.Range("K74").Value = 5810.778
$ws.Range("M74").Value = -4874.778
# Row 77
$ws.Range("H77").Value = 6561.8184
$ws.Range("I77").Value = 5810.778
$ws.Range("K77").Value = 29053.89
$ws.Range("M77").Value = -24373.89
# Row 100
$ws.Range("H100").Value = 43123.543
$ws.Range("I100").Value = 46571.137
$ws.Range("J100").Value = 5200
$ws.Range("K100").Value = 46571.137
$ws.Range("L100").Value = 5200
$ws.Range("M100").Value = -46030.137
$ws.Range("N100").Value = -6282
# Row 106
$ws.Range("H106").Value = 23123.414
$ws.Range("I106").Value = 13088.667
$ws.Range("K106").Value = 13088.667
$ws.Range("M106").Value = -12457.667
# Row 113
$ws.Range("H113").Value = 4644.7617
$ws.Range("I113").Value = 3995.3845
$ws.Range("J113").Value = 5700
$ws.Range("K113").Value = 3995.3845
$ws.Range("L113").Value = 5700
$ws.Range("M113").Value = -741.3845000000001
$ws.Range("N113").Value = -12208
# Row 116
$ws.Range("H116").Value = 11261.096
$ws.Range("I116").Value = 2603.6667
$ws.Range("J116").Value = 17754.166
$ws.Range("K116").Value = 2603.6667
$ws.Range("L116").Value = 17754.166
$ws.Range("M116").Value = 838.3332999999998
$ws.Range("N116").Value = -24638.166
# Row 138
$ws.Range("H138").Value = 4810.9185
$ws.Range("J138").Value = 5958.7744
$ws.Range("L138").Value = 17876.3232
$ws.Range("N138").Value = -28156.3232

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4895.5713
$ws.Range("I2").Value = 4476.148
$ws.Range("J2").Value = 6311.125
$ws.Range("K2").Value = 4476.148
$ws.Range("L2").Value = 6311.125
$ws.Range("M2").Value = -4363.148
$ws.Range("N2").Value = -6537.125
# Row 31
$ws.Range("H31").Value = 11616.75
$ws.Range("I31").Value = 11616.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 11616.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -11322.75
$ws.Range("N31").ClearContents()
# Row 32
$ws.Range("H32").Value = 5119.9316
$ws.Range("I32").Value = 5119.9316
$ws.Range("K32").Value = 5119.9316
$ws.Range("M32").Value = -4832.9316
# Row 45
$ws.Range("H45").Value = 10022.923
$ws.Range("I45").Value = 11407.091
$ws.Range("J45").Value = 2410
$ws.Range("K45").Value = 11407.091
$ws.Range("L45").Value = 2410
$ws.Range("M45").Value = -11030.091
$ws.Range("N45").Value = -3164
# Row 63
$ws.Range("H63").Value = 3261.75
$ws.Range("I63").Value = 1717.6
$ws.Range("K63").Value = 1717.6
$ws.Range("M63").Value = -1031.6
# Row 66
$ws.Range("H66").Value = 3261.75
$ws.Range("I66").Value = 1717.6
$ws.Range("K66").Value = 8588
$ws.Range("M66").Value = -5156
# Row 102
$ws.Range("H102").Value = 1457.381
$ws.Range("I102").Value = 1537.9412
$ws.Range("J102").Value = 1115
$ws.Range("K102").Value = 1537.9412
$ws.Range("L102").Value = 1115
$ws.Range("M102").Value = 84.05880000000002
$ws.Range("N102").Value = -4359
# Row 116
$ws.Range("H116").Value = 4895.5713
$ws.Range("I116").Value = 4476.148
$ws.Range("J116").Value = 6311.125
$ws.Range("K116").Value = 4476.148
$ws.Range("L116").Value = 6311.125
$ws.Range("M116").Value = -2182.148
$ws.Range("N116").Value = -10899.125
# Row 132
$ws.Range("H132").Value = 3598.6667
$ws.Range("I132").Value = 3668.6
$ws.Range("K132").Value = 11005.8
$ws.Range("M132").Value = -8475.799999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4895.5713
$ws.Range("I3").Value = 4476.148
$ws.Range("J3").Value = 6311.125
$ws.Range("K3").Value = 4476.148
$ws.Range("L3").Value = 6311.125
$ws.Range("M3").Value = -4362.148
$ws.Range("N3").Value = -6539.125
# Row 22
$ws.Range("H22").Value = 1248.75
$ws.Range("I22").Value = 998.3333
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 998.3333
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -825.3333
$ws.Range("N22").Value = -2346
# Row 26
$ws.Range("H26").Value = 8635
$ws.Range("I26").Value = 7602
$ws.Range("J26").Value = 13800
$ws.Range("K26").Value = 7602
$ws.Range("L26").Value = 13800
$ws.Range("M26").Value = -7310
$ws.Range("N26").Value = -14384
# Row 80
$ws.Range("H80").Value = 706.5333000000001
$ws.Range("I80").Value = 176.85715
$ws.Range("J80").Value = 1170
$ws.Range("K80").Value = 176.85715
$ws.Range("L80").Value = 1170
$ws.Range("M80").Value = 821.14285
$ws.Range("N80").Value = -3166
# Row 82
$ws.Range("H82").Value = 51322.555
$ws.Range("I82").Value = 11539.333
$ws.Range("K82").Value = 11539.333
$ws.Range("M82").Value = -11156.333
# Row 83
$ws.Range("H83").Value = 706.5333000000001
$ws.Range("I83").Value = 176.85715
$ws.Range("J83").Value = 1170
$ws.Range("K83").Value = 884.28575
$ws.Range("L83").Value = 5850
$ws.Range("M83").Value = 4107.71425
$ws.Range("N83").Value = -15834
# Row 85
$ws.Range("H85").Value = 51322.555
$ws.Range("I85").Value = 11539.333
$ws.Range("K85").Value = 11539.333
$ws.Range("M85").Value = -10213.333
# Row 86
$ws.Range("H86").Value = 591971.6
$ws.Range("I86").Value = 1668698.4
$ws.Range("J86").Value = 4666.1816
$ws.Range("K86").Value = 1668698.4
$ws.Range("L86").Value = 4666.1816
$ws.Range("M86").Value = -1667575.4
$ws.Range("N86").Value = -6912.1816
# Row 89
$ws.Range("H89").Value = 591971.6
$ws.Range("I89").Value = 1668698.4
$ws.Range("J89").Value = 4666.1816
$ws.Range("K89").Value = 8343492
$ws.Range("L89").Value = 23330.908
$ws.Range("M89").Value = -8337876
$ws.Range("N89").Value = -34562.908
# Row 94
$ws.Range("H94").Value = 3355.75
$ws.Range("I94").Value = 3719.4
$ws.Range("J94").Value = 2749.6667
$ws.Range("K94").Value = 3719.4
$ws.Range("L94").Value = 2749.6667
$ws.Range("M94").Value = -3268.4
$ws.Range("N94").Value = -3651.6667
# Row 102
$ws.Range("H102").Value = 17123.715
$ws.Range("I102").Value = 17123.715
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 17123.715
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -13878.715
$ws.Range("N102").ClearContents()
# Row 105
$ws.Range("H105").Value = 3974.4
$ws.Range("I105").Value = 3555.8462
$ws.Range("K105").Value = 3555.8462
$ws.Range("M105").Value = -1808.8462
# Row 107
$ws.Range("H107").Value = 4869
$ws.Range("I107").Value = 4935.1816
$ws.Range("K107").Value = 4935.1816
$ws.Range("M107").Value = -3015.1816

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 9702.333000000001
$ws.Range("I16").Value = 9702.333000000001
$ws.Range("K16").Value = 9702.333000000001
$ws.Range("M16").Value = -9415.333000000001
# Row 86
$ws.Range("H86").Value = 3432.625
$ws.Range("I86").Value = 3098.6667
$ws.Range("J86").Value = 3766.5833
$ws.Range("K86").Value = 3098.6667
$ws.Range("L86").Value = 3766.5833
$ws.Range("M86").Value = -1975.6667
$ws.Range("N86").Value = -6012.5833
# Row 89
$ws.Range("H89").Value = 3432.625
$ws.Range("I89").Value = 3098.6667
$ws.Range("J89").Value = 3766.5833
$ws.Range("K89").Value = 15493.3335
$ws.Range("L89").Value = 18832.9165
$ws.Range("M89").Value = -9877.333500000001
$ws.Range("N89").Value = -30064.9165
# Row 99
$ws.Range("H99").Value = 9163.554
$ws.Range("I99").Value = 5282.0586
$ws.Range("J99").Value = 10855.487
$ws.Range("K99").Value = 5282.0586
$ws.Range("L99").Value = 10855.487
$ws.Range("M99").Value = -3784.0586
$ws.Range("N99").Value = -13851.487
# Row 105
$ws.Range("H105").Value = 710.4706
$ws.Range("I105").Value = 642.4375
$ws.Range("J105").Value = 1799
$ws.Range("K105").Value = 642.4375
$ws.Range("L105").Value = 1799
$ws.Range("M105").Value = 1104.5625
$ws.Range("N105").Value = -5293
# Row 107
$ws.Range("H107").Value = 1831.2
$ws.Range("I107").Value = 1925.1428
$ws.Range("J107").Value = 1612
$ws.Range("K107").Value = 1925.1428
$ws.Range("L107").Value = 1612
$ws.Range("M107").Value = -5.142800000000079
$ws.Range("N107").Value = -5452
# Row 113
$ws.Range("H113").Value = 9702.333000000001
$ws.Range("I113").Value = 9702.333000000001
$ws.Range("K113").Value = 9702.333000000001
$ws.Range("M113").Value = -7532.333000000001
# Row 126
$ws.Range("H126").Value = 9163.554
$ws.Range("I126").Value = 5282.0586
$ws.Range("J126").Value = 10855.487
$ws.Range("K126").Value = 15846.1758
$ws.Range("L126").Value = 32566.461
$ws.Range("M126").Value = -13376.1758
$ws.Range("N126").Value = -37506.461
# Row 132
$ws.Range("H132").Value = 28781.188
$ws.Range("J132").Value = 44888.05
$ws.Range("L132").Value = 134664.15
$ws.Range("N132").Value = -139724.15
# Row 134
$ws.Range("H134").Value = 10289
$ws.Range("I134").Value = 5579.2
$ws.Range("K134").Value = 16737.6
$ws.Range("M134").Value = -14202.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 12000
$ws.Range("I58").Value = 12000
$ws.Range("K58").Value = 36000
$ws.Range("M58").Value = -35872
# Row 107
$ws.Range("H107").Value = 2295.6667
$ws.Range("I107").Value = 2386.75
$ws.Range("J107").Value = 2222.8
$ws.Range("K107").Value = 7160.25
$ws.Range("L107").Value = 6668.400000000001
$ws.Range("M107").Value = -5240.25
$ws.Range("N107").Value = -10508.4
# Row 113
$ws.Range("H113").Value = 1908.091
$ws.Range("I113").Value = 2078.1667
$ws.Range("K113").Value = 6234.500100000001
$ws.Range("M113").Value = -4064.500100000001
# Row 116
$ws.Range("H116").Value = 5601
$ws.Range("I116").Value = 7575
$ws.Range("J116").Value = 3627
$ws.Range("K116").Value = 22725
$ws.Range("L116").Value = 10881
$ws.Range("M116").Value = -19283
$ws.Range("N116").Value = -17765
# Row 122
$ws.Range("H122").Value = 14286521
$ws.Range("I122").Value = 749
$ws.Range("J122").Value = 16667483
$ws.Range("K122").Value = 6741
$ws.Range("L122").Value = 150007347
$ws.Range("M122").Value = -4291
$ws.Range("N122").Value = -150012247
# Row 131
$ws.Range("H131").Value = 3051.7334
$ws.Range("I131").Value = 1400
$ws.Range("J131").Value = 3554.4348
$ws.Range("K131").Value = 4200
$ws.Range("L131").Value = 10663.3044
$ws.Range("M131").Value = 840
$ws.Range("N131").Value = -20743.3044
# Row 132
$ws.Range("H132").Value = 3706283.8
$ws.Range("I132").Value = 2785.2144
$ws.Range("J132").Value = 7694666.5
$ws.Range("K132").Value = 25066.9296
$ws.Range("L132").Value = 69251998.5
$ws.Range("M132").Value = -22536.9296
$ws.Range("N132").Value = -69257058.5
# Row 140
$ws.Range("H140").Value = 2096.3333
$ws.Range("I140").Value = 827.2222
$ws.Range("K140").Value = 2481.6666
$ws.Range("M140").Value = 2698.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 10833.167
$ws.Range("I29").Value = 4500
$ws.Range("J29").Value = 13999.75
$ws.Range("K29").Value = 4500
$ws.Range("L29").Value = 13999.75
$ws.Range("M29").Value = -4210
$ws.Range("N29").Value = -14579.75
# Row 58
$ws.Range("H58").Value = 35333
$ws.Range("J58").Value = 38000
$ws.Range("L58").Value = 38000
$ws.Range("N58").Value = -38554
# Row 62
$ws.Range("H62").Value = 500085
$ws.Range("J62").Value = 500085
$ws.Range("L62").Value = 500085
$ws.Range("N62").Value = -501457
# Row 65
$ws.Range("H65").Value = 500085
$ws.Range("J65").Value = 500085
$ws.Range("L65").Value = 1500255
$ws.Range("N65").Value = -1507119
# Row 97
$ws.Range("H97").Value = 1917.3334
$ws.Range("I97").Value = 1854
$ws.Range("J97").Value = 1949
$ws.Range("K97").Value = 1854
$ws.Range("L97").Value = 1949
$ws.Range("M97").Value = -1358
$ws.Range("N97").Value = -2941
# Row 107
$ws.Range("H107").Value = 291.82352
$ws.Range("I107").Value = 290.53333
$ws.Range("J107").Value = 301.5
$ws.Range("K107").Value = 290.53333
$ws.Range("L107").Value = 301.5
$ws.Range("M107").Value = 1629.46667
$ws.Range("N107").Value = -4141.5
# Row 126
$ws.Range("H126").Value = 4525.933
$ws.Range("I126").Value = 3987.5557
$ws.Range("J126").Value = 5333.5
$ws.Range("K126").Value = 11962.6671
$ws.Range("L126").Value = 16000.5
$ws.Range("M126").Value = -9492.667099999999
$ws.Range("N126").Value = -20940.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5334.2856
$ws.Range("I7").Value = 5225.1665
$ws.Range("J7").Value = 5989
$ws.Range("K7").Value = 5225.1665
$ws.Range("L7").Value = 5989
$ws.Range("M7").Value = -5113.1665
$ws.Range("N7").Value = -6213
# Row 46
$ws.Range("H46").Value = 1415.1915
$ws.Range("I46").Value = 961.9286
$ws.Range("K46").Value = 961.9286
$ws.Range("M46").Value = -773.9286
# Row 74
$ws.Range("H74").Value = 38333
$ws.Range("I74").Value = 22500
$ws.Range("K74").Value = 22500
$ws.Range("M74").Value = -21502
# Row 77
$ws.Range("H77").Value = 38333
$ws.Range("I77").Value = 22500
$ws.Range("K77").Value = 67500
$ws.Range("M77").Value = -62508
# Row 93
$ws.Range("H93").Value = 4839.276
$ws.Range("I93").Value = 5131.4585
$ws.Range("J93").Value = 3436.8
$ws.Range("K93").Value = 5131.4585
$ws.Range("L93").Value = 3436.8
$ws.Range("M93").Value = -3883.4585
$ws.Range("N93").Value = -5932.8
# Row 98
$ws.Range("H98").Value = 30118.334
$ws.Range("J98").Value = 30118.334
$ws.Range("L98").Value = 30118.334
$ws.Range("N98").Value = -36108.334
# Row 99
$ws.Range("H99").Value = 40884.9
$ws.Range("I99").Value = 40884.9
$ws.Range("K99").Value = 40884.9
$ws.Range("M99").Value = -37889.9
# Row 122
$ws.Range("H122").Value = 5500.32
$ws.Range("I122").Value = 4652.6924
$ws.Range("K122").Value = 13958.0772
$ws.Range("M122").Value = -11508.0772
# Row 126
$ws.Range("H126").Value = 5334.2856
$ws.Range("I126").Value = 5225.1665
$ws.Range("J126").Value = 5989
$ws.Range("K126").Value = 15675.4995
$ws.Range("L126").Value = 17967
$ws.Range("M126").Value = -13205.4995
$ws.Range("N126").Value = -22907

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 4999
$ws.Range("J5").Value = 4999
$ws.Range("L5").Value = 4999
$ws.Range("N5").Value = -5223
# Row 21
$ws.Range("H21").Value = 24998.666
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19765
# Row 35
$ws.Range("H35").Value = 24998.666
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19710
# Row 61
$ws.Range("H61").Value = 40367.125
$ws.Range("I61").Value = 40367.125
$ws.Range("K61").Value = 40367.125
$ws.Range("M61").Value = -40075.125
# Row 70
$ws.Range("H70").Value = 67244.5
$ws.Range("J70").Value = 69992.664
$ws.Range("L70").Value = 69992.664
$ws.Range("N70").Value = -70622.664
# Row 73
$ws.Range("H73").Value = 67244.5
$ws.Range("J73").Value = 69992.664
$ws.Range("L73").Value = 69992.664
$ws.Range("N73").Value = -72176.664
# Row 75
$ws.Range("H75").Value = 34950.5
$ws.Range("I75").Value = 34950.5
$ws.Range("K75").Value = 34950.5
$ws.Range("M75").Value = -34014.5
# Row 78
$ws.Range("H78").Value = 34950.5
$ws.Range("I78").Value = 34950.5
$ws.Range("K78").Value = 104851.5
$ws.Range("M78").Value = -100171.5
# Row 107
$ws.Range("H107").Value = 2398.8333
$ws.Range("I107").Value = 1798.3334
$ws.Range("J107").Value = 2999.3333
$ws.Range("K107").Value = 5395.0002
$ws.Range("L107").Value = 8997.999899999999
$ws.Range("M107").Value = -3475.0002
$ws.Range("N107").Value = -12837.9999
# Row 122
$ws.Range("H122").Value = 3856.2354
$ws.Range("I122").Value = 3050.818
$ws.Range("J122").Value = 5332.8335
$ws.Range("K122").Value = 9152.454000000002
$ws.Range("L122").Value = 15998.5005
$ws.Range("M122").Value = -6702.454000000002
$ws.Range("N122").Value = -20898.5005
# Row 132
$ws.Range("H132").Value = 116585.13
$ws.Range("I132").Value = 137377.11
$ws.Range("J132").Value = 35496.4
$ws.Range("K132").Value = 412131.33
$ws.Range("L132").Value = 106489.2
$ws.Range("M132").Value = -409601.33
$ws.Range("N132").Value = -111549.2
# Row 136
$ws.Range("H136").Value = 9377204
$ws.Range("I136").Value = 13638333
$ws.Range("K136").Value = 40914999
$ws.Range("M136").Value = -40912449

